# Weekly update: insert a new daily record at row 49 (Coliflor, Vega
# Monumental Concepción), pushing the existing rows (49-164) down by one
# row. This grows the used range from A1:R164 to A1:R165.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 49..164 down to 50..165, leaving a fresh blank row 49.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new price record.
$ws.Cells.Item(49, 1).Value  = 11
$ws.Cells.Item(49, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value  = "Bíobío"
$ws.Cells.Item(49, 4).Value  = 44498
$ws.Cells.Item(49, 5).Value  = 8
$ws.Cells.Item(49, 6).Value  = 100112008
$ws.Cells.Item(49, 7).Value  = "Coliflor"
$ws.Cells.Item(49, 8).Value  = "Sin especificar"
$ws.Cells.Item(49, 9).Value  = "Primera"
$ws.Cells.Item(49, 10).Value = 2500
$ws.Cells.Item(49, 11).Value = 600
$ws.Cells.Item(49, 12).Value = 650
$ws.Cells.Item(49, 13).Value = 630
$ws.Cells.Item(49, 14).Value = "`$/unidad"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 630
$ws.Cells.Item(49, 17).Value = 1
$ws.Cells.Item(49, 18).Value = "Hortaliza"

# Apply the same date number format used by the rest of column D.
$ws.Cells.Item(49, 4).NumberFormat = $ws.Cells.Item(50, 4).NumberFormat()
